$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ G=0.5979736666666666; H=1.793921; I=0.03342655292740804; J=0.03342655292740804;
            M=0.668273; N=2.004819; O=0.01328414746766746; P=0.01328414746766746;
            Q=0.3996096561443333; R=3.596486905298999; S=0.0004440432584234799; T=0.0004440432584234798 }
    3  = @{ G=0.5979736666666666; H=1.793921; I=0.03342655292740804; J=0.03342655292740804;
            O=0.3831531055114357; P=0.3831531055114357;
            Q=11.52589438778211; R=103.733049490039; S=0.01280748756067876; T=0.01280748756067876 }
    4  = @{ G=0.5979736666666666; H=1.793921; I=0.03342655292740804; J=0.03342655292740804;
            M=30.36285833333334; N=91.08857500000001; O=0.6035627470208969; P=0.6035627470208967;
            Q=18.15618972806389; R=163.405707552575; S=0.0201750221083058; T=0.0201750221083058 }
    5  = @{ I=0.8874158839838097; J=0.8874158839838097;
            M=0.668273; N=2.004819; O=0.01328414746766746; P=0.01328414746766746;
            Q=10.60892988355433; R=95.48036895198899; S=0.01178856346799141; T=0.0117885634679914 }
    6  = @{ I=0.8874158839838097; J=0.8874158839838097;
            O=0.3831531055114357; P=0.3831531055114357;
            S=0.3400161518285726; T=0.3400161518285726 }
    7  = @{ I=0.8874158839838097; J=0.8874158839838097;
            M=30.36285833333334; N=91.08857500000001; O=0.6035627470208969; P=0.6035627470208967;
            R=4338.132643550825; S=0.5356111686872457; T=0.5356111686872456 }
    8  = @{ I=0.07915756308878232; J=0.07915756308878232;
            M=0.668273; N=2.004819; O=0.01328414746766746; P=0.01328414746766746;
            Q=0.9463173374720002; R=8.516856037248001; S=0.001051540741252575; T=0.001051540741252574 }
    9  = @{ I=0.07915756308878232; J=0.07915756308878232;
            O=0.3831531055114357; P=0.3831531055114357;
            Q=27.29451984285867; S=0.03032946612218434; T=0.03032946612218433 }
    10 = @{ I=0.07915756308878232; J=0.07915756308878232;
            M=30.36285833333334; N=91.08857500000001; O=0.6035627470208969; P=0.6035627470208967;
            Q=42.99575062293334; R=386.9617556064001; S=0.0477765562253454; T=0.04777655622534539 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
